$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 3.4
$ws.Range("I4").Value = 2.2
$ws.Range("T4").Value = 8
$ws.Range("W4").Value = 41
$ws.Range("AA4").Value = 6.5
$ws.Range("G6").Value = 2.02
$ws.Range("H6").Value = 2.95
$ws.Range("I6").Value = 3.9
$ws.Range("M6").Value = 2.45
$ws.Range("N6").Value = 2.25
$ws.Range("O6").Value = 1.5
$ws.Range("R6").Value = 1.93
$ws.Range("S6").Value = 1.7
$ws.Range("U6").Value = 8.5
$ws.Range("V6").Value = 8.75
$ws.Range("W6").Value = 18.5
$ws.Range("X6").Value = 19
$ws.Range("Y6").Value = 37
$ws.Range("Z6").Value = 6.9
$ws.Range("AA6").Value = 5.8
$ws.Range("AB6").Value = 16.5
$ws.Range("AC6").Value = 100
$ws.Range("AE6").Value = 9.25
$ws.Range("AF6").Value = 20
$ws.Range("AG6").Value = 13.5
$ws.Range("AH6").Value = 65
$ws.Range("AI6").Value = 40
$ws.Range("AJ6").Value = 55
$ws.Range("N10").Value = 1.83
$ws.Range("O10").Value = 2.03
$ws.Range("G12").Value = 2.52
$ws.Range("H12").Value = 3.25
$ws.Range("I12").Value = 2.45
$ws.Range("L12").Value = 1.26
$ws.Range("M12").Value = 3.5
$ws.Range("N12").Value = 1.72
$ws.Range("O12").Value = 1.88
$ws.Range("R12").Value = 1.66
$ws.Range("S12").Value = 2.09
$ws.Range("T12").Value = 8.25
$ws.Range("U12").Value = 11.5
$ws.Range("V12").Value = 8
$ws.Range("W12").Value = 23
$ws.Range("X12").Value = 16
$ws.Range("Y12").Value = 20
$ws.Range("Z12").Value = 11
$ws.Range("AA12").Value = 5.7
$ws.Range("AB12").Value = 10.25
$ws.Range("AC12").Value = 37
$ws.Range("AD12").Value = 200
$ws.Range("AE12").Value = 7.9
$ws.Range("AF12").Value = 11
$ws.Range("AG12").Value = 8
$ws.Range("AH12").Value = 22
$ws.Range("AI12").Value = 16
$ws.Range("AJ12").Value = 21
$ws.Range("P13").Value = 1.22
$ws.Range("P14").Value = 1.36
$ws.Range("G15").Value = 2.3
$ws.Range("H15").Value = 2.9
$ws.Range("I15").Value = 3.45
$ws.Range("J15").Value = 1.12
$ws.Range("K15").Value = 5.6
$ws.Range("L15").Value = 1.53
$ws.Range("M15").Value = 2.42
$ws.Range("N15").Value = 2.55
$ws.Range("O15").Value = 1.47
$ws.Range("P15").Value = 1.57
$ws.Range("Q15").Value = 2.32
$ws.Range("R15").Value = 2.05
$ws.Range("S15").Value = 1.7
$ws.Range("U15").Value = 10.5
$ws.Range("V15").Value = 10.25
$ws.Range("W15").Value = 26
$ws.Range("Z15").Value = 5.6
$ws.Range("AA15").Value = 6
$ws.Range("AB15").Value = 19.5
$ws.Range("AE15").Value = 7.4
$ws.Range("AF15").Value = 17.5
$ws.Range("AG15").Value = 13.5
$ws.Range("AI15").Value = 45
$ws.Range("AJ15").Value = 65
$ws.Range("G16").Value = 2.85
$ws.Range("H16").Value = 2.85
$ws.Range("I16").Value = 2.72
$ws.Range("K16").Value = 5.6
$ws.Range("Q16").Value = 2.37
$ws.Range("T16").Value = 6.8
$ws.Range("U16").Value = 14
$ws.Range("V16").Value = 11.5
$ws.Range("W16").Value = 40
$ws.Range("X16").Value = 32
$ws.Range("Z16").Value = 5.6
$ws.Range("AA16").Value = 5.8
$ws.Range("AB16").Value = 17.5
$ws.Range("AC16").Value = 110
$ws.Range("AE16").Value = 6.8
$ws.Range("AF16").Value = 13.5
$ws.Range("AG16").Value = 11
$ws.Range("AH16").Value = 37
$ws.Range("AI16").Value = 30
$ws.Range("G23").Value = 3.65
$ws.Range("H23").Value = 3.7
$ws.Range("I23").Value = 1.88
$ws.Range("J23").Value = 1.03
$ws.Range("K23").Value = 9
$ws.Range("L23").Value = 1.17
$ws.Range("M23").Value = 4.4
$ws.Range("N23").Value = 1.53
$ws.Range("O23").Value = 2.32
$ws.Range("P23").Value = 1.29
$ws.Range("Q23").Value = 3.3
$ws.Range("R23").Value = 1.5
$ws.Range("S23").Value = 2.42
$ws.Range("T23").Value = 16
$ws.Range("U23").Value = 25
$ws.Range("W23").Value = 55
$ws.Range("X23").Value = 27
$ws.Range("Y23").Value = 26
$ws.Range("Z23").Value = 9
$ws.Range("AA23").Value = 7.6
$ws.Range("AB23").Value = 11.5
$ws.Range("AC23").Value = 37
$ws.Range("AD23").Value = 200
$ws.Range("AE23").Value = 10
$ws.Range("AF23").Value = 11
$ws.Range("AI23").Value = 13
$ws.Range("AJ23").Value = 18.5
